$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1774193548387097
$ws.Range("C2").Value = 0.6064516129032258
$ws.Range("J2").Value = 0.01935483870967742
$ws.Range("P2").Value = 0.1354838709677419
$ws.Range("S2").Value = 0.06129032258064516

$ws.Range("B3").Value = 0.005025125628140704
$ws.Range("C3").Value = 0.03015075376884422
$ws.Range("J3").Value = 0.03015075376884422
$ws.Range("P3").Value = 0.7889447236180904
$ws.Range("S3").Value = 0.1457286432160804

$ws.Range("J4").Value = 0.05172413793103448
$ws.Range("P4").Value = 0.6724137931034483
$ws.Range("S4").Value = 0.2758620689655172

$ws.Range("J5").Value = 0.25
$ws.Range("S5").Value = 0.25

$ws.Range("B6").Value = 0.09900990099009901
$ws.Range("D6").Value = 0.03465346534653466
$ws.Range("F6").Value = 0.08415841584158416
$ws.Range("J6").Value = 0.2475247524752475
$ws.Range("O6").Value = 0.0198019801980198
$ws.Range("Q6").Value = 0.1831683168316832
$ws.Range("R6").Value = 0.08415841584158416
$ws.Range("S6").Value = 0.2475247524752475

$ws.Range("B7").Value = 0.1442307692307692
$ws.Range("D7").Value = 0.03365384615384615
$ws.Range("F7").Value = 0.07692307692307693
$ws.Range("J7").Value = 0.1394230769230769
$ws.Range("O7").Value = 0.01923076923076923
$ws.Range("Q7").Value = 0.2067307692307692
$ws.Range("R7").Value = 0.04807692307692308
$ws.Range("S7").Value = 0.3317307692307692

$ws.Range("B8").Value = 0.08746355685131195
$ws.Range("D8").Value = 0.02623906705539359
$ws.Range("E8").Value = 0.002915451895043732
$ws.Range("F8").Value = 0.04664723032069971
$ws.Range("J8").Value = 0.1049562682215743
$ws.Range("O8").Value = 0.01166180758017493
$ws.Range("Q8").Value = 0.1749271137026239
$ws.Range("R8").Value = 0.09329446064139942
$ws.Range("S8").Value = 0.4518950437317784

$ws.Range("B9").Value = 0.09090909090909091
$ws.Range("D9").Value = 0.006493506493506494
$ws.Range("E9").Value = 0.006493506493506494
$ws.Range("F9").Value = 0.07142857142857142
$ws.Range("J9").Value = 0.1428571428571428
$ws.Range("O9").Value = 0.01948051948051948
$ws.Range("Q9").Value = 0.1753246753246753
$ws.Range("R9").Value = 0.08441558441558442
$ws.Range("S9").Value = 0.4025974025974026

$ws.Range("B10").Value = 0.1312803889789303
$ws.Range("D10").Value = 0.02836304700162074
$ws.Range("E10").Value = 0.002431118314424636
$ws.Range("F10").Value = 0.05672609400324149
$ws.Range("J10").Value = 0.1329011345218801
$ws.Range("O10").Value = 0.01215559157212318
$ws.Range("Q10").Value = 0.2204213938411669
$ws.Range("R10").Value = 0.08670988654781199
$ws.Range("S10").Value = 0.3290113452188007

$ws.Range("F11").Value = 0.003115264797507788
$ws.Range("G11").Value = 0.1308411214953271
$ws.Range("J11").Value = 0.08411214953271028
$ws.Range("K11").Value = 0.1838006230529595
$ws.Range("L11").Value = 0.5825545171339563
$ws.Range("S11").Value = 0.01557632398753894

$ws.Range("G12").Value = 0.7604166666666666
$ws.Range("J12").Value = 0.171875
$ws.Range("K12").Value = 0.01041666666666667
$ws.Range("L12").Value = 0.02083333333333333
$ws.Range("S12").Value = 0.03645833333333334

$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3111111111111111
$ws.Range("S13").Value = 0.02222222222222222

$ws.Range("F15").Value = 0.04761904761904762
$ws.Range("H15").Value = 0.164021164021164
$ws.Range("I15").Value = 0.06878306878306878
$ws.Range("J15").Value = 0.3650793650793651
$ws.Range("K15").Value = 0.06349206349206349
$ws.Range("M15").Value = 0.01587301587301587
$ws.Range("O15").Value = 0.1216931216931217
$ws.Range("S15").Value = 0.1534391534391534

$ws.Range("F16").Value = 0.008733624454148471
$ws.Range("H16").Value = 0.1266375545851528
$ws.Range("I16").Value = 0.04803493449781659
$ws.Range("J16").Value = 0.5152838427947598
$ws.Range("K16").Value = 0.1222707423580786
$ws.Range("M16").Value = 0.01746724890829694
$ws.Range("O16").Value = 0.03056768558951965
$ws.Range("S16").Value = 0.1310043668122271

$ws.Range("F17").Value = 0.02017937219730942
$ws.Range("H17").Value = 0.1345291479820628
$ws.Range("I17").Value = 0.08295964125560538
$ws.Range("J17").Value = 0.4641255605381166
$ws.Range("K17").Value = 0.08968609865470852
$ws.Range("M17").Value = 0.02017937219730942
$ws.Range("O17").Value = 0.06502242152466367
$ws.Range("S17").Value = 0.1233183856502242

$ws.Range("F18").Value = 0.00546448087431694
$ws.Range("H18").Value = 0.1366120218579235
$ws.Range("I18").Value = 0.08743169398907104
$ws.Range("J18").Value = 0.3825136612021858
$ws.Range("K18").Value = 0.1092896174863388
$ws.Range("M18").Value = 0.01092896174863388
$ws.Range("O18").Value = 0.09836065573770492
$ws.Range("S18").Value = 0.1693989071038251

$ws.Range("F19").Value = 0.0196078431372549
$ws.Range("H19").Value = 0.1755793226381462
$ws.Range("I19").Value = 0.0659536541889483
$ws.Range("J19").Value = 0.3538324420677362
$ws.Range("K19").Value = 0.1381461675579323
$ws.Range("M19").Value = 0.02584670231729055
$ws.Range("O19").Value = 0.05793226381461675
$ws.Range("S19").Value = 0.1631016042780749

